$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row at position 3 for "Korean Air" (pushes Asiana/HongKong/Rossya/Peach down)
$ws.Rows.Item(3).Insert()

# 2. Fill in the new row 3 data (airline + airplane first; status text set further below
#    so that new shared strings get appended in the same order as the authored workbook)
$ws.Range("A3").Value = "Korean Air"
$ws.Range("B3").Value = "A320-251N (A321N)"

# 3. Add column E header "Last update", matching the header style of row 1
$ws.Range("D1").Copy() | Out-Null
$ws.Range("E1").PasteSpecial(-4122) | Out-Null
$ws.Range("E1").Value = "Last update"

# 4. Row 3 status text
$ws.Range("D3").Value = "Engine done + main colors"

# 5. Update row 2 (Uzbekistan Airways) status text
$ws.Range("D2").Value = "Completed"

# 6. Column E width
$ws.Columns.Item(5).ColumnWidth = 13.6

# 7. Dates for rows 2 and 3 ("last update" column), centered + short-date number format
$ws.Range("E2").HorizontalAlignment = -4108
$ws.Range("E2").VerticalAlignment = -4108
$ws.Range("E2").NumberFormat = "mm-dd-yy"
$ws.Range("E2").Value = 45317

$ws.Range("E3").HorizontalAlignment = -4108
$ws.Range("E3").VerticalAlignment = -4108
$ws.Range("E3").NumberFormat = "mm-dd-yy"
$ws.Range("E3").Value = 45317

# 8. Blank centered cells for E4:E11 (rest of the "last update" column placeholders)
$ws.Range("E4:E11").HorizontalAlignment = -4108
$ws.Range("E4:E11").VerticalAlignment = -4108

# 9. Status colors in column C
$ws.Range("C2").Interior.Color = 5287936   # green  FF00B050 - Completed
$ws.Range("C3").Interior.Color = 65535     # yellow FFFFFF00 - Engine done + main colors
$ws.Range("C4:C7").Interior.Color = 192    # red    FFC00000 - Not started

# 10. View state
$excel.ActiveWindow.Zoom = 174
$ws.Range("D4").Select() | Out-Null
